# daily auto push: 2026-02-16 03:18 UTC
# Insert a new daily-ranking entry for 2026/02/16 07:xx ahead of the
# 2026/12/29 block, shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 819 (and everything below it) down one row, leaving a blank
# row 819 for the new entry.
$ws.Rows(819).Insert()

# Column A holds the date as literal text (e.g. "2026/12/29"), not a real
# Excel date serial. Force text interpretation before writing the value so
# it isn't auto-parsed into a date, then drop the temporary number format
# so the cell keeps the sheet's default (unstyled) formatting.
$ws.Range("A819").NumberFormat = "@"
$ws.Range("A819").Value = "2026/02/16"
$ws.Range("A819").ClearFormats()

$ws.Range("B819").Value = "月"
$ws.Range("C819").Value = 7
$ws.Range("D819").Value = 201
